# Apply schema renames described in the commit diff:
#  - Person sheet:   A1 "name"   -> "last_name"
#  - Author sheet:   F1 "name"   -> "last_name"
#  - ImageSize sheet: A1 "height" -> "height_im", B1 "width" -> "width_im"

$wb = $excel.ActiveWorkbook

$wsPerson = $wb.Worksheets.Item("Person")
$wsPerson.Range("A1").Value = "last_name"

$wsAuthor = $wb.Worksheets.Item("Author")
$wsAuthor.Range("F1").Value = "last_name"

$wsImageSize = $wb.Worksheets.Item("ImageSize")
$wsImageSize.Range("A1").Value = "height_im"
$wsImageSize.Range("B1").Value = "width_im"
